$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format so numeric-looking strings
# (e.g. "1.002") are stored as inline/shared strings rather than numbers,
# matching the source workbook where these are plain text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.050.08"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.820.69"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "233.83"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "0.5979"
$ws.Range("E6").Value = "  -4.31%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.06922"
$ws.Range("E8").Value = "  -6.09%  "
$ws.Range("D9").Value = "0.2740"
$ws.Range("E9").Value = "  -4.94%  "
$ws.Range("D10").Value = "23.16"
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("D11").Value = "0.07582"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "1.837.66"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "4.722"
$ws.Range("E13").Value = "  -4.59%  "
$ws.Range("D14").Value = "0.6220"
$ws.Range("E14").Value = "  -6.16%  "
$ws.Range("D15").Value = "0.000009656"
$ws.Range("E15").Value = "  -8.22%  "
$ws.Range("D16").Value = "77.13"
$ws.Range("E16").Value = "  -5.07%  "
$ws.Range("D17").Value = "28.677.50"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "5.557"
$ws.Range("E18").Value = "  -11.00%  "
$ws.Range("D19").Value = "215.12"
$ws.Range("E19").Value = "  -7.76%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "11.47"
$ws.Range("E21").Value = "  -6.15%  "
$ws.Range("E22").Value = "  -6.19%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "155.98"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "7.890"
$ws.Range("E25").Value = "  -6.31%  "
$ws.Range("E26").Value = "  -4.31%  "
$ws.Range("D27").Value = "16.37"
$ws.Range("E27").Value = "  -5.04%  "
$ws.Range("D28").Value = "0.06374"
$ws.Range("E28").Value = "  -10.58%  "
$ws.Range("D29").Value = "1.417"
$ws.Range("E29").Value = "  -4.70%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").Value = "3.822"
$ws.Range("E31").Value = "  -4.97%  "
$ws.Range("D32").Value = "3.732"
$ws.Range("E32").Value = "  -7.38%  "
$ws.Range("D33").Value = "1.716"
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("D34").Value = "1.085"
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("D35").Value = "0.6426"
$ws.Range("E35").Value = "  -7.64%  "
$ws.Range("D36").Value = "2.534"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").Value = "2.735"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").Value = "0.01745"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "6.520"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "1.140.58"
$ws.Range("E40").Value = "  -7.43%  "
$ws.Range("D41").Value = "0.8790"
$ws.Range("E41").Value = "  -7.15%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "1.973.86"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "100.04"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "61.34"
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").Value = "0.00000000113"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "1.598"
$ws.Range("E47").Value = "  -4.78%  "
$ws.Range("E48").Value = "  -5.63%  "
$ws.Range("D49").Value = "0.05503"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").Value = "0.4523"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "6.397"
$ws.Range("E51").Value = "  -7.74%  "

# Revert the Price column back to its original (default) style now that
# the values are committed as text, so no stray cell-style index remains.
$ws.Range("D2:D51").Style = "Normal"

